$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 326 (shifts old rows 326-384 down to 327-385,
# carrying their formatting/styles with them - matches Excel's native Insert behaviour).
$ws.Rows.Item(326).Insert()

# Populate the newly inserted row 326 with the new weekly data point.
$ws.Cells.Item(326, 1).Value = 10
$ws.Cells.Item(326, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(326, 3).Value = "La Araucanía"
$ws.Cells.Item(326, 4).Value = 44694
$ws.Cells.Item(326, 5).Value = 9
$ws.Cells.Item(326, 6).Value = 100114014
$ws.Cells.Item(326, 7).Value = "Betarraga"
$ws.Cells.Item(326, 8).Value = "Sin especificar"
$ws.Cells.Item(326, 9).Value = "Primera"
$ws.Cells.Item(326, 10).Value = 45
$ws.Cells.Item(326, 11).Value = 9000
$ws.Cells.Item(326, 12).Value = 9000
$ws.Cells.Item(326, 13).Value = 9000
$ws.Cells.Item(326, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(326, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(326, 16).Value = 750
$ws.Cells.Item(326, 17).Value = 12
$ws.Cells.Item(326, 18).Value = "Hortaliza"
